# Insert a new data row right before current row 38. This shifts all the
# existing rows 38:113 down to 39:114 (pushing the previous last record,
# old row 113, out to the new row 114), and the newly inserted row 38 is
# populated with a fresh price record for Tuna at "Vega Modelo de Temuco".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above row 38 (keeps formatting from the row above,
# same as a normal Excel "Insert" on a selected row).
$ws.Rows.Item(38).Insert()

# Populate the newly inserted row 38 with the new record's data.
$ws.Cells.Item(38, 1).Value  = 10
$ws.Cells.Item(38, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(38, 3).Value  = "La Araucanía"
$ws.Cells.Item(38, 4).Value  = "08/24/2023"
$ws.Cells.Item(38, 5).Value  = 9
$ws.Cells.Item(38, 6).Value  = "Fruta"
$ws.Cells.Item(38, 7).Value  = 100107
$ws.Cells.Item(38, 8).Value  = "Otros"
$ws.Cells.Item(38, 9).Value  = 100107011
$ws.Cells.Item(38, 10).Value = "Tuna"
$ws.Cells.Item(38, 11).Value = "Sin especificar"
$ws.Cells.Item(38, 12).Value = "Especial"
$ws.Cells.Item(38, 13).Value = 120
$ws.Cells.Item(38, 14).Value = 33000
$ws.Cells.Item(38, 15).Value = 33000
$ws.Cells.Item(38, 16).Value = 33000
$ws.Cells.Item(38, 17).Value = "`$/caja 16 kilos"
$ws.Cells.Item(38, 18).Value = "Provincia de Los Andes"
$ws.Cells.Item(38, 19).Value = 2062
$ws.Cells.Item(38, 20).Value = 16
